# updated batching of commands
# Extends the LINUX-COMMANDS sheet with two new command rows ("who" /
# "which terminal" and "ls; who" / "batching commands"), repeats the
# header row across columns C:L, restyles the populated cells with a
# monospaced "Courier New" font, widens columns A/B, refreshes the
# header/footer caption text, and moves the active selection to A6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Repeat the header row (LINUX COMMAND / MEANING) across C1:L1 ---
$ws.Range("C1").Value = "LINUX COMMAND"
$ws.Range("D1").Value = "MEANING"
$ws.Range("E1").Value = "LINUX COMMAND"
$ws.Range("F1").Value = "MEANING"
$ws.Range("G1").Value = "LINUX COMMAND"
$ws.Range("H1").Value = "MEANING"
$ws.Range("I1").Value = "LINUX COMMAND"
$ws.Range("J1").Value = "MEANING"
$ws.Range("K1").Value = "LINUX COMMAND"
$ws.Range("L1").Value = "MEANING"

# --- New data rows: batching commands ---
$ws.Range("A4").Value = "who"
$ws.Range("B4").Value = "which terminal"
$ws.Range("A5").Value = "ls; who"
$ws.Range("B5").Value = "batching commands"

# --- Restyle every populated cell with Courier New 11 ---
$fontRanges = @("A1:L1", "A2", "A3", "A4:B4", "A5:B5")
foreach ($addr in $fontRanges) {
    $fnt = $ws.Range($addr).Font
    $fnt.Name = "Courier New"
    $fnt.Size = 11
}

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 19.7
$ws.Columns.Item(2).ColumnWidth = 26.48

# --- Header / footer caption text ---
$pageSetup = $ws.PageSetup
$pageSetup.CenterHeader = '&"Times New Roman,Normál"&12&A'
$pageSetup.CenterFooter = '&"Times New Roman,Normál"&12Page &P'

# --- Move active selection to A6 ---
$null = $ws.Range("A6").Select()
